# The sheet gains one new data row: a new row 244 is inserted (shifting the
# former rows 244:301 down to 245:302), and the new row 244 is populated
# with a fresh "Choclero" price record for "Región del Maule".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("244:244").Insert()

$ws.Range("A244").Value = 7
$ws.Range("B244").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C244").Value = 'Ñuble'
$ws.Range("D244").Value = 45005
$ws.Range("E244").Value = 16
$ws.Range("F244").Value = 100112024
$ws.Range("G244").Value = 'Choclo'
$ws.Range("H244").Value = 'Choclero'
$ws.Range("I244").Value = 'Primera'
$ws.Range("J244").Value = 5000
$ws.Range("K244").Value = 400
$ws.Range("L244").Value = 400
$ws.Range("M244").Value = 400
$ws.Range("N244").Value = '$/unidad'
$ws.Range("O244").Value = 'Región del Maule'
$ws.Range("P244").Value = 400
$ws.Range("Q244").Value = 1
$ws.Range("R244").Value = 'Hortaliza'
